# Adds two new weekly price rows (Primera/Segunda) for Albahaca at
# "Mercado Mayorista Lo Valledor de Santiago" dated 2022-01-17 (serial 44578),
# inserted right before the current row 330, pushing all following rows down
# by two and growing the used range from A1:R357 to A1:R359.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block that starts at row 330.
$ws.Rows("330:331").Insert()

# New row 330: Calidad "Primera"
$ws.Range("A330").Value = 6
$ws.Range("B330").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C330").Value = "Metropolitana"
$ws.Range("D330").Value = 44578
$ws.Range("E330").Value = 13
$ws.Range("F330").Value = 100112052
$ws.Range("G330").Value = "Albahaca"
$ws.Range("H330").Value = "Sin especificar"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 490
$ws.Range("K330").Value = 3000
$ws.Range("L330").Value = 4500
$ws.Range("M330").Value = 3765
$ws.Range("N330").Value = "`$/docena de matas"
$ws.Range("O330").Value = "Región Metropolitana"
$ws.Range("P330").Value = 628
$ws.Range("Q330").Value = 6
$ws.Range("R330").Value = "Hortaliza"

# New row 331: Calidad "Segunda"
$ws.Range("A331").Value = 6
$ws.Range("B331").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C331").Value = "Metropolitana"
$ws.Range("D331").Value = 44578
$ws.Range("E331").Value = 13
$ws.Range("F331").Value = 100112052
$ws.Range("G331").Value = "Albahaca"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Segunda"
$ws.Range("J331").Value = 220
$ws.Range("K331").Value = 2500
$ws.Range("L331").Value = 3500
$ws.Range("M331").Value = 3091
$ws.Range("N331").Value = "`$/docena de matas"
$ws.Range("O331").Value = "Región Metropolitana"
$ws.Range("P331").Value = 515
$ws.Range("Q331").Value = 6
$ws.Range("R331").Value = "Hortaliza"

# Make sure the date cells keep the same date/time number format used by
# every other "Fecha" cell in column D.
$ws.Range("D330").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D331").NumberFormat = "YYYY-MM-DD HH:MM:SS"
